# Generate Report for Handoff
# Adds a new tracked file (c2ddd944-...) as row 3 on the Overview, zh-cn and
# de-de sheets, mirroring the existing a920dcca-... row that is already
# present as row 2 on every sheet.

$wb = $excel.ActiveWorkbook

$newFileName = 'c2ddd944-a457-435b-a5f4-a88fd1c4ad13ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$newPath     = 'e2e\c2ddd944-a457-435b-a5f4-a88fd1c4ad13ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$statusReady = 'Ready for handoff'
$dateOverview = '2016-09-02 18:30:24'
$zhXlf       = 'c2ddd944-a457-435b-a5f4-a88fd1c4ad13oooooooooooooooooooooooooooooooooooooooo.7d4d07005e65c6cacf76eb400a7295c2b094bbae.zh-cn.xlf'
$zhDate      = '2016-09-02 18:30:09'
$deXlf       = 'c2ddd944-a457-435b-a5f4-a88fd1c4ad13oooooooooooooooooooooooooooooooooooooooo.7d4d07005e65c6cacf76eb400a7295c2b094bbae.de-de.xlf'
$deDate      = '2016-09-02 18:30:24'

$hyperlinkTarget = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68265ffbf3230686bf200c9e5f1cfa40da366ce8/e2e/' + $newFileName

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(3,1).Value = $newFileName
$wsOverview.Cells.Item(3,2).Value = $newPath
$wsOverview.Cells.Item(3,3).Value = ".md"
$wsOverview.Range("D2").Copy($wsOverview.Range("D3"))
$wsOverview.Cells.Item(3,5).Value = $statusReady
$wsOverview.Cells.Item(3,6).Value = $statusReady
$wsOverview.Cells.Item(3,7).Value = $dateOverview
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, "", "", $newPath) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Cells.Item(3,1).Value = $newFileName
$wsZh.Cells.Item(3,2).Value = ".md"
$wsZh.Cells.Item(3,3).Value = $statusReady
$wsZh.Cells.Item(3,4).Value = "e2e"
$wsZh.Cells.Item(3,5).Value = "ht"
$wsZh.Cells.Item(3,6).Value = "False"
$wsZh.Cells.Item(3,7).Value = $zhXlf
$wsZh.Cells.Item(3,8).Value = $zhDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("I2").Copy($wsZh.Range("I3"))
$wsZh.Range("J2").Copy($wsZh.Range("J3"))

$wsZh.Cells.Item(3,11).Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("L2").Copy($wsZh.Range("L3"))

$wsZh.Cells.Item(3,13).Value = "True"

$wsZh.Range("N2").Copy($wsZh.Range("N3"))

$wsZh.Cells.Item(3,15).Value = "False"

$wsZh.Range("P2").Copy($wsZh.Range("P3"))

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hyperlinkTarget, "", "", $newFileName) | Out-Null
$wsZh.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Cells.Item(3,1).Value = $newFileName
$wsDe.Cells.Item(3,2).Value = ".md"
$wsDe.Cells.Item(3,3).Value = $statusReady
$wsDe.Cells.Item(3,4).Value = "e2e"
$wsDe.Cells.Item(3,5).Value = "ht"
$wsDe.Cells.Item(3,6).Value = "False"
$wsDe.Cells.Item(3,7).Value = $deXlf
$wsDe.Cells.Item(3,8).Value = $deDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("I2").Copy($wsDe.Range("I3"))
$wsDe.Range("J2").Copy($wsDe.Range("J3"))

$wsDe.Cells.Item(3,11).Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("L2").Copy($wsDe.Range("L3"))

$wsDe.Cells.Item(3,13).Value = "True"

$wsDe.Range("N2").Copy($wsDe.Range("N3"))

$wsDe.Cells.Item(3,15).Value = "False"

$wsDe.Range("P2").Copy($wsDe.Range("P3"))

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hyperlinkTarget, "", "", $newFileName) | Out-Null
$wsDe.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Column width tweaks (E/F on Overview, C on zh-cn/de-de) widen slightly
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZh.Columns.Item(3).ColumnWidth = 16.3
$wsDe.Columns.Item(3).ColumnWidth = 16.3
